$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '50.776.44'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -1.02%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.923.53'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -1.80%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '374.67'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  -1.83%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '99.90'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -2.37%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.535'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -1.60%  '

$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.583'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -1.17%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '35.91'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  -2.67%  '

$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -0.56%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0841'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -0.04%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '3.380.32'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -1.81%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '17.96'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  -1.09%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.54'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +0.43%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '11.35'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +54.11%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.912.42'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -2.21%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.985'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -0.86%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '50.674.22'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -1.04%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '3.05'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -6.02%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '12.32'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -3.63%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.0₃0950'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -0.62%  '

$ws.Cells.Item(23, 2).NumberFormat = '@'
$ws.Cells.Item(23, 2).Value = 'BitcoinCash'
$ws.Cells.Item(23, 3).NumberFormat = '@'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '264.58'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +1.69%  '

$ws.Cells.Item(24, 2).NumberFormat = '@'
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).NumberFormat = '@'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '68.68'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -0.32%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.14'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +8.52%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '7.95'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -1.85%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.34'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -1.53%  '

$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +0.05%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '25.42'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -1.73%  '

$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -3.74%  '

$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -7.42%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '9.95'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +1.20%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '50.75'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  -0.21%  '

$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -0.85%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '32.95'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -3.98%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0436'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -3.63%  '

$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -0.02%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.08'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  +3.60%  '

$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -0.29%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '16.39'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -3.35%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.80'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -1.00%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.45'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -4.19%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '119.39'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -2.83%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '20.88'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -2.55%  '

$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -1.91%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.35'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +2.69%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.33'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -1.34%  '

$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -2.42%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.981.15'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -2.33%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0324'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -2.12%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '5.17'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +1.31%  '
